$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "RubyOnRails"
$ws.Range("C5").Value = "http://stackoverflow.com/questions/13222406/undefined-method-current-page-for-array0x007fd5ef6dd158-kaminari"
$ws.Range("B5").Value = "undifind method 'current_page' for array kaminari gem."

$ws.Range("B5").Select()
